$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 799.5
$ws.Range("I12").Value = 799.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 799.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -629.5
$ws.Range("N12").ClearContents()

# Row 18
$ws.Range("H18").Value = 373.14285
$ws.Range("I18").Value = 204
$ws.Range("K18").Value = 204
$ws.Range("M18").Value = 80

# Row 74
$ws.Range("H74").Value = 14289858
$ws.Range("I74").Value = 25003500
$ws.Range("K74").Value = 25003500
$ws.Range("M74").Value = -25002564

# Row 77
$ws.Range("H77").Value = 14289858
$ws.Range("I77").Value = 25003500
$ws.Range("K77").Value = 125017500
$ws.Range("M77").Value = -125012820

# Row 113
$ws.Range("H113").Value = 5071.0586
$ws.Range("I113").Value = 3001.6667
$ws.Range("J113").Value = 6199.8184
$ws.Range("K113").Value = 3001.6667
$ws.Range("L113").Value = 6199.8184
$ws.Range("M113").Value = 252.3332999999998
$ws.Range("N113").Value = -12707.8184

# Row 116
$ws.Range("H116").Value = 232530.66
$ws.Range("I116").Value = 456856
$ws.Range("J116").Value = 8205.317999999999
$ws.Range("K116").Value = 456856
$ws.Range("L116").Value = 8205.317999999999
$ws.Range("M116").Value = -453414
$ws.Range("N116").Value = -15089.318

# Row 123
$ws.Range("H123").Value = 41807.5
$ws.Range("J123").Value = 41807.5
$ws.Range("L123").Value = 41807.5
$ws.Range("N123").Value = -51607.5

# Row 132
$ws.Range("H132").Value = 184712.88
$ws.Range("I132").Value = 2908.25
$ws.Range("J132").Value = 1431373.1
$ws.Range("K132").Value = 8724.75
$ws.Range("L132").Value = 4294119.300000001
$ws.Range("M132").Value = -6194.75
$ws.Range("N132").Value = -4299179.300000001

# Row 138
$ws.Range("H138").Value = 3398.25
$ws.Range("I138").Value = 1056
$ws.Range("J138").Value = 4058.8845
$ws.Range("K138").Value = 3168
$ws.Range("L138").Value = 12176.6535
$ws.Range("M138").Value = 1972
$ws.Range("N138").Value = -22456.6535

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6069.456
$ws.Range("I32").Value = 5144.878
$ws.Range("J32").Value = 8438.6875
$ws.Range("K32").Value = 5144.878
$ws.Range("L32").Value = 8438.6875
$ws.Range("M32").Value = -4857.878
$ws.Range("N32").Value = -9012.6875

# Row 63
$ws.Range("H63").Value = 9237661
$ws.Range("I63").Value = 15392613
$ws.Range("J63").Value = 5233.3335
$ws.Range("K63").Value = 15392613
$ws.Range("L63").Value = 5233.3335
$ws.Range("M63").Value = -15391927
$ws.Range("N63").Value = -6605.3335

# Row 64
$ws.Range("H64").Value = 34250
$ws.Range("J64").Value = 34250
$ws.Range("L64").Value = 34250
$ws.Range("N64").Value = -34746

# Row 66
$ws.Range("H66").Value = 9237661
$ws.Range("I66").Value = 15392613
$ws.Range("J66").Value = 5233.3335
$ws.Range("K66").Value = 76963065
$ws.Range("L66").Value = 26166.6675
$ws.Range("M66").Value = -76959633
$ws.Range("N66").Value = -33030.6675

# Row 67
$ws.Range("H67").Value = 34250
$ws.Range("J67").Value = 34250
$ws.Range("L67").Value = 34250
$ws.Range("N67").Value = -35966

# Row 132
$ws.Range("H132").Value = 2661.2307
$ws.Range("I132").Value = 1434.7273
$ws.Range("J132").Value = 9407
$ws.Range("K132").Value = 4304.1819
$ws.Range("L132").Value = 28221
$ws.Range("M132").Value = -1774.1819
$ws.Range("N132").Value = -33281

# Row 133
$ws.Range("H133").Value = 38260
$ws.Range("J133").Value = 38260
$ws.Range("L133").Value = 38260
$ws.Range("N133").Value = -43320

# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# Row 135
$ws.Range("H135").Value = 48489.855
$ws.Range("J135").Value = 48489.855
$ws.Range("L135").Value = 48489.855
$ws.Range("N135").Value = -58629.855

$ws = $wb.Worksheets.Item("BSM")
# Row 132
$ws.Range("H132").Value = 56902.223
$ws.Range("J132").Value = 56902.223
$ws.Range("L132").Value = 56902.223
$ws.Range("N132").Value = -67022.223

$ws = $wb.Worksheets.Item("CRP")
# Row 123
$ws.Range("H123").Value = 40835
$ws.Range("J123").Value = 40835
$ws.Range("L123").Value = 40835
$ws.Range("N123").Value = -50635

# Row 127
$ws.Range("H127").Value = 41868.332
$ws.Range("J127").Value = 41868.332
$ws.Range("L127").Value = 41868.332
$ws.Range("N127").Value = -51788.332

# Row 132
$ws.Range("H132").Value = 1814.9714
$ws.Range("I132").Value = 1258.7587
$ws.Range("K132").Value = 3776.2761
$ws.Range("M132").Value = -1246.2761

# Row 133
$ws.Range("H133").Value = 70333.336
$ws.Range("J133").Value = 70333.336
$ws.Range("L133").Value = 70333.336
$ws.Range("N133").Value = -75393.336

$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 350
$ws.Range("J97").Value = 350
$ws.Range("L97").Value = 1050
$ws.Range("N97").Value = -2042

# Row 113
$ws.Range("H113").Value = 549.12
$ws.Range("I113").Value = 528.6316
$ws.Range("J113").Value = 614
$ws.Range("K113").Value = 1585.8948
$ws.Range("L113").Value = 1842
$ws.Range("M113").Value = 584.1052
$ws.Range("N113").Value = -6182

# Row 131
$ws.Range("H131").Value = 13514399
$ws.Range("I131").Value = 83333736
$ws.Range("J131").Value = 978.9355
$ws.Range("K131").Value = 250001208
$ws.Range("L131").Value = 2936.8065
$ws.Range("M131").Value = -249996168
$ws.Range("N131").Value = -13016.8065

$ws = $wb.Worksheets.Item("LTW")
# Row 69
$ws.Range("H69").Value = 40000
$ws.Range("J69").Value = 40000
$ws.Range("L69").Value = 40000
$ws.Range("N69").Value = -41622

# Row 72
$ws.Range("H72").Value = 40000
$ws.Range("J72").Value = 40000
$ws.Range("L72").Value = 120000
$ws.Range("N72").Value = -128112

# Row 122
$ws.Range("H122").Value = 3042.4443
$ws.Range("I122").Value = 1854.8334
$ws.Range("K122").Value = 5564.5002
$ws.Range("M122").Value = -3114.5002

# Row 139
$ws.Range("H139").Value = 43096.25
$ws.Range("J139").Value = 43096.25
$ws.Range("L139").Value = 43096.25
$ws.Range("N139").Value = -53376.25

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 10102718
$ws.Range("I132").Value = 1102.6086
$ws.Range("J132").Value = 33336434
$ws.Range("K132").Value = 3307.8258
$ws.Range("L132").Value = 100009302
$ws.Range("M132").Value = -777.8258000000001
$ws.Range("N132").Value = -100014362
